# "arreglos en los ingresos y costos"
#
# Updates the volume figures on "Volumenes Ingresos" (row 6 & 7), which
# ripple through "Ingresos", "Ganancias" and "Indicadores Financieros" via
# formulas. Also refreshes the cached external-reference costs on
# "Ganancias" (B3/C3, sourced from the external "Gastos Resultado"
# workbook) and replays the various cell selections / active-sheet the
# author ended up on.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Volumenes Ingresos": corrected volumes for Semestre 1/2 - 2017
# ---------------------------------------------------------------------
$volumenes = $wb.Worksheets.Item("Volumenes Ingresos")
$volumenes.Range("B6").Value = 23
$volumenes.Range("C6").Value = 7
$volumenes.Range("B7").Value = 21
$volumenes.Range("C7").Value = 9

# ---------------------------------------------------------------------
# 2) "Ganancias": refresh the externally-linked cost figures
#    ('[1]Gastos Resultado'!$B$4 / $C$4). The external workbook isn't
#    reachable from here, so write the refreshed cached numbers directly.
# ---------------------------------------------------------------------
$ganancias = $wb.Worksheets.Item("Ganancias")
$ganancias.Range("B3").Value = 4794714.41
$ganancias.Range("C3").Value = 8309827.3409200003

# ---------------------------------------------------------------------
# 3) Replay selections left on each sheet. Do the sheets that should end
#    up NOT active first; the final Select() below decides the workbook's
#    active tab.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Tipo Ingreso").Range("C3").Select() | Out-Null
$volumenes.Range("D6").Select() | Out-Null
$wb.Worksheets.Item("Ingresos").Range("C11").Select() | Out-Null

# "Indicadores Financieros" becomes the active sheet/tab (matches
# activeTab going from 2 -> 4, and tabSelected moving off "Ingresos").
$indicadores = $wb.Worksheets.Item("Indicadores Financieros")
$indicadores.Columns.Item(3).ColumnWidth = 12.86
$indicadores.Range("B4").Select() | Out-Null
